# Thermometer calibration feature: adds a "confirm" UI string near the top
# of the language table, plus a "thermometer_calibration_result" / "target"
# pair of rows right after the existing "thermometer_calibration" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 8: confirm / CONFIRM (right after "close"/"CLOSE") ---
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).Value2 = "confirm"
$ws.Cells.Item(8, 2).Value2 = "CONFIRM"

# --- Insert two new rows (17-18) right after "thermometer_calibration" ---
$ws.Range("A17:A18").EntireRow.Insert()
$ws.Cells.Item(17, 1).Value2 = "thermometer_calibration_result"
$ws.Cells.Item(17, 2).Value2 = "Thermometer Calibration Result"
$ws.Cells.Item(18, 1).Value2 = "target"
$ws.Cells.Item(18, 2).Value2 = "Target"

# --- Update the on-screen view to follow the newly inserted rows ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A18").Select()
